$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row 2 (A2:H2) to upper case column names
$ws.Range("A2").Value = "STATION"
$ws.Range("B2").Value = "NAME"
$ws.Range("C2").Value = "NETID"
$ws.Range("D2").Value = "SOURCE"
$ws.Range("E2").Value = "LAT"
$ws.Range("F2").Value = "LON"
$ws.Range("G2").Value = "DISTANCE"
$ws.Range("H2").Value = "INTENSITY"

# Widen column H slightly (character-width units; this renders as width="10" in the
# underlying OOXML column definition, which previously was width="7.5")
$ws.Columns.Item(8).ColumnWidth = 9.17

# Update the selected range in the sheet view
$ws.Range("I2:M2").Select()
